$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.768.60'
$ws.Range("E2").Value = '  -0.91%  '
$ws.Range("D3").Value = '1.928.05'
$ws.Range("E3").Value = '  -1.43%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.39'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.37%  '
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4860'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.46%  '
$ws.Range("E8").Value = '  -0.58%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06833'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("E10").Value = '  -0.27%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '106.28'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.94%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07768'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.57%  '
$ws.Range("D13").Value = '1.928.46'
$ws.Range("E13").Value = '  -1.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.327'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.52%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6964'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.75%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '274.75'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -3.11%  '
$ws.Range("D17").Value = '30.757.01'
$ws.Range("E17").Value = '  -1.00%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007673'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.46%  '
$ws.Range("E19").Value = '  +0.01%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.584'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.72%  '
$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.98'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.78%  '
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.452'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.62%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.852'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '164.98'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -2.87%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '19.45'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.79%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.146'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.46%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1038'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.90%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.387'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.74%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.563'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.551'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.20%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.350'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04878'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.32%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7570'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.29%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.142'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.65%  '
$ws.Range("E36").Value = '  -0.02%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.721'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01987'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.645'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.21%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.470'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.96%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '77.51'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +3.54%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.056'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.95%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8836'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.17%  '
$ws.Range("E44").Value = '  -0.98%  '
$ws.Range("E45").Value = '  -1.67%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.899'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -3.58%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.001'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '982.00'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.88%  '
$ws.Range("E49").Value = '  -1.85%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.11'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.98%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.228'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.12%  '
